# DepositsCalculationScheme.xlsx - "deposits view shows forecast"
#
# The sheet "Deposit estimations" mini class-diagram on Лист2 is updated:
#  - the old free-floating "DepositCalculationData" textbox (a stale class
#    dump) is removed from the drawing canvas;
#  - the "DepositEstimations" field list (column F, rows 36-41) is replaced
#    with the new forecast-oriented fields;
#  - the "DepositRateLine" code block (column K, rows 15-20) is moved left
#    to column H so it lines up under the now-wider column F;
#  - the active selection / column F width are refreshed to match the new
#    layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# 1. Remove the stale "DepositCalculationData" textbox shape.
$ws.Shapes.Item("TextBox 1").Delete()

# 2. Move the "DepositRateLine" code block from column K to column H.
for ($r = 15; $r -le 20; $r++) {
    $ws.Cells.Item($r, 8).Value2 = $ws.Cells.Item($r, 11).Value2
}
$ws.Range("H15").Font.Bold = $true
$ws.Range("K15:K20").Clear()

# 3. Update the "DepositEstimations" field list to show the forecast fields.
$ws.Range("F36").Value2 = "Period PeriodForThisMonthPayment"
$ws.Range("F38").Value2 = "Period PeriodForUpToEndPayment"
$ws.Range("F37").Value2 = "decimal ProcentsInThisMonth"
$ws.Range("F39").Value2 = "decimal ProcentsUpToFinish"
$ws.Range("F40").Value2 = "decimal DevaluationInUsd"
$ws.Range("F41").Value2 = "decimal ProfitInUsd"

# 4. Widen column F so the longer field names fit, and refresh the selection.
$ws.Columns.Item(6).ColumnWidth = 32.877604166666664
$ws.Range("I25").Select()
